# Fix errors in column names, as per
# https://github.com/rfordatascience/tidytuesday/issues/396
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6: score_team2 description was mistakenly "team one score" -> "Team two score"
$ws.Range("C6").Value = "Team two score"

# Row 7: field name was "wickets_team2" but it actually documents team one's wickets
$ws.Range("A7").Value = "wickets_team1"

# Row 8: field name was "wickets_team" (missing the "2") and its description
# incorrectly described team one; both now correctly reference team two
$ws.Range("A8").Value = "wickets_team2"
$ws.Range("C8").Value = "wickets fallen for team two; if 10 it means all out. "

# Update the selected cell to match the saved view state
$ws.Range("C8").Select()
